# Auto-generated edit script applying the numeric corrections described in the
# commit diff (per-leve currentAveragePrice / LevePrice / LeveProfit recalculations).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1480.5385
$ws.Range("J17").Value = 1480.5385
$ws.Range("L17").Value = 4441.6155
$ws.Range("N17").Value = -4777.6155

$ws.Range("H43").Value = 1750
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 1750
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 1750
$ws.Range("M43").Value = $null
$ws.Range("N43").Value = -1888

$ws.Range("H52").Value = 3000
$ws.Range("J52").Value = 3000
$ws.Range("L52").Value = 9000
$ws.Range("N52").Value = -9320

$ws.Range("H113").Value = 8645.182000000001
$ws.Range("J113").Value = 9071.143
$ws.Range("L113").Value = 9071.143
$ws.Range("N113").Value = -15579.143

$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = $null
$ws.Range("N138").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4098.5
$ws.Range("I61").Value = 4148.5
$ws.Range("J61").Value = 3998.5
$ws.Range("K61").Value = 4148.5
$ws.Range("L61").Value = 3998.5
$ws.Range("M61").Value = -3936.5
$ws.Range("N61").Value = -4422.5

$ws.Range("H119").Value = 38333
$ws.Range("J119").Value = 38333
$ws.Range("L119").Value = 38333
$ws.Range("N119").Value = -48009

$ws.Range("H136").Value = 4098.5
$ws.Range("I136").Value = 4148.5
$ws.Range("J136").Value = 3998.5
$ws.Range("K136").Value = 12445.5
$ws.Range("L136").Value = 11995.5
$ws.Range("M136").Value = -9895.5
$ws.Range("N136").Value = -17095.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4067.5
$ws.Range("I86").Value = 4501
$ws.Range("J86").Value = 1900
$ws.Range("K86").Value = 4501
$ws.Range("L86").Value = 1900
$ws.Range("M86").Value = -3378
$ws.Range("N86").Value = -4146

$ws.Range("H89").Value = 4067.5
$ws.Range("I89").Value = 4501
$ws.Range("J89").Value = 1900
$ws.Range("K89").Value = 22505
$ws.Range("L89").Value = 9500
$ws.Range("M89").Value = -16889
$ws.Range("N89").Value = -20732

$ws.Range("H109").Value = 79990
$ws.Range("J109").Value = 79990
$ws.Range("L109").Value = 79990
$ws.Range("N109").Value = -82764

$ws.Range("H110").Value = 100000
$ws.Range("J110").Value = 100000
$ws.Range("L110").Value = 100000
$ws.Range("N110").Value = -108180

$ws.Range("H134").Value = 6037.364
$ws.Range("I134").Value = 6037.364
$ws.Range("K134").Value = 18112.092
$ws.Range("M134").Value = -15577.092

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 540.6667
$ws.Range("I16").Value = 540.6667
$ws.Range("K16").Value = 540.6667
$ws.Range("M16").Value = -253.6667

$ws.Range("H31").Value = 4293.4287
$ws.Range("I31").Value = 2831.6667
$ws.Range("J31").Value = 5389.75
$ws.Range("K31").Value = 2831.6667
$ws.Range("L31").Value = 5389.75
$ws.Range("M31").Value = -2536.6667
$ws.Range("N31").Value = -5979.75

$ws.Range("H34").Value = 4293.4287
$ws.Range("I34").Value = 2831.6667
$ws.Range("J34").Value = 5389.75
$ws.Range("K34").Value = 2831.6667
$ws.Range("L34").Value = 5389.75
$ws.Range("M34").Value = -2629.6667
$ws.Range("N34").Value = -5793.75

$ws.Range("H113").Value = 540.6667
$ws.Range("I113").Value = 540.6667
$ws.Range("K113").Value = 540.6667
$ws.Range("M113").Value = 1629.3333

$ws.Range("H132").Value = 3997
$ws.Range("I132").Value = 3996
$ws.Range("K132").Value = 11988
$ws.Range("M132").Value = -9458

$ws.Range("H134").Value = 1046.091
$ws.Range("I134").Value = 1023
$ws.Range("K134").Value = 3069
$ws.Range("M134").Value = -534

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").Value = $null

$ws.Range("H32").Value = 995
$ws.Range("J32").Value = 995
$ws.Range("L32").Value = 2985
$ws.Range("N32").Value = -3551

$ws.Range("I131").Value = 1123.25
$ws.Range("J131").Value = 1999.5
$ws.Range("K131").Value = 3369.75
$ws.Range("L131").Value = 5998.5
$ws.Range("M131").Value = 1670.25
$ws.Range("N131").Value = -16078.5

$ws.Range("H134").Value = 500999.5
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2276
$ws.Range("J80").Value = 2476.5
$ws.Range("L80").Value = 2476.5
$ws.Range("N80").Value = -4472.5

$ws.Range("H83").Value = 2276
$ws.Range("J83").Value = 2476.5
$ws.Range("L83").Value = 12382.5
$ws.Range("N83").Value = -22366.5

$ws.Range("H97").Value = 1568.4546
$ws.Range("I97").Value = 450.33334
$ws.Range("K97").Value = 450.33334
$ws.Range("M97").Value = 45.66665999999998

$ws.Range("H132").Value = 1666.3334
$ws.Range("I132").Value = 1249.5
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 3748.5
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -1218.5
$ws.Range("N132").Value = -12560

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 929.1111
$ws.Range("I9").Value = 881.4286
$ws.Range("K9").Value = 881.4286
$ws.Range("M9").Value = -657.4286

$ws.Range("H46").Value = 2940.476
$ws.Range("I46").Value = 2654
$ws.Range("K46").Value = 2654
$ws.Range("M46").Value = -2466

$ws.Range("H122").Value = 3252
$ws.Range("I122").Value = 3252
$ws.Range("K122").Value = 9756
$ws.Range("M122").Value = -7306

$ws.Range("H132").Value = 4868.6665
$ws.Range("I132").Value = 4868.6665
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 14605.9995
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -12075.9995
$ws.Range("N132").Value = $null

$ws.Range("H136").Value = 4206.143
$ws.Range("I136").Value = 2221.5
$ws.Range("K136").Value = 6664.5
$ws.Range("M136").Value = -4114.5
